$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to be treated as text so numeric-looking
# strings (e.g. "226.86", "0.0521") are not auto-converted to numbers,
# matching the original inline-string cell content.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '33.922.79'
$ws.Range('E2').Value = '  +10.90%  '

$ws.Range('D3').Value = '1.806.01'
$ws.Range('E3').Value = '  +7.46%  '

$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.47%  '

$ws.Range('D5').Value = '226.86'
$ws.Range('E5').Value = '  +3.21%  '

$ws.Range('D6').Value = '0.538'
$ws.Range('E6').Value = '  +2.45%  '

$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.38%  '

$ws.Range('D8').Value = '31.15'
$ws.Range('E8').Value = '  +3.21%  '

$ws.Range('D9').Value = '46.95'
$ws.Range('E9').Value = '  +5.88%  '

$ws.Range('E10').Value = '  +5.80%  '

$ws.Range('D11').Value = '0.0663'
$ws.Range('E11').Value = '  +5.94%  '

$ws.Range('D12').Value = '0.0925'
$ws.Range('E12').Value = '  +1.91%  '

$ws.Range('D13').Value = '2.063.74'
$ws.Range('E13').Value = '  +7.05%  '

$ws.Range('D14').Value = '1.812.13'
$ws.Range('E14').Value = '  +7.56%  '

$ws.Range('D15').Value = '0.634'
$ws.Range('E15').Value = '  +2.53%  '

$ws.Range('D16').Value = '33.853.36'
$ws.Range('E16').Value = '  +10.45%  '

$ws.Range('D17').Value = '10.08'
$ws.Range('E17').Value = '  -3.11%  '

$ws.Range('D18').Value = '4.24'
$ws.Range('E18').Value = '  +6.78%  '

$ws.Range('D19').Value = '69.05'
$ws.Range('E19').Value = '  +4.23%  '

$ws.Range('D20').Value = '255.25'
$ws.Range('E20').Value = '  +4.38%  '

$ws.Range('D21').Value = '0.0₃0741'
$ws.Range('E21').Value = '  +3.90%  '

$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  -0.04%  '

$ws.Range('D23').Value = '10.41'
$ws.Range('E23').Value = '  +2.38%  '

$ws.Range('D24').Value = '4.30'
$ws.Range('E24').Value = '  +0.79%  '

$ws.Range('E25').Value = '  +0.54%  '

$ws.Range('D26').Value = '157.25'
$ws.Range('E26').Value = '  -0.36%  '

$ws.Range('D27').Value = '2.14'
$ws.Range('E27').Value = '  +423.82%  '

$ws.Range('D28').Value = '16.45'
$ws.Range('E28').Value = '  +3.85%  '

$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D29').Value = '7.04'
$ws.Range('E29').Value = '  +5.39%  '

$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D30').Value = '0.114'
$ws.Range('E30').Value = '  +3.31%  '

$ws.Range('E31').Value = '  -0.35%  '

$ws.Range('D32').Value = '3.82'
$ws.Range('E32').Value = '  +10.34%  '

$ws.Range('D33').Value = '0.0508'
$ws.Range('E33').Value = '  +2.07%  '

$ws.Range('D34').Value = '1.19'
$ws.Range('E34').Value = '  +4.87%  '

$ws.Range('D35').Value = '3.48'
$ws.Range('E35').Value = '  +6.25%  '

$ws.Range('D36').Value = '1.532.10'
$ws.Range('E36').Value = '  +1.44%  '

$ws.Range('E37').Value = '  +2.51%  '

$ws.Range('D38').Value = '1.06'
$ws.Range('E38').Value = '  +2.73%  '

$ws.Range('D39').Value = '83.65'
$ws.Range('E39').Value = '  -0.42%  '

$ws.Range('D40').Value = '0.0185'
$ws.Range('E40').Value = '  +4.10%  '

$ws.Range('D41').Value = '0.613'
$ws.Range('E41').Value = '  +4.62%  '

$ws.Range('E43').Value = '  +1.54%  '

$ws.Range('D44').Value = '0.902'
$ws.Range('E44').Value = '  +7.97%  '

$ws.Range('E45').Value = '  +6.65%  '

$ws.Range('D46').Value = '0.0521'
$ws.Range('E46').Value = '  +4.08%  '

$ws.Range('E47').Value = '  +4.08%  '

$ws.Range('D48').Value = '1.952.41'
$ws.Range('E48').Value = '  +7.01%  '

$ws.Range('D49').Value = '0.997'
$ws.Range('E49').Value = '  -0.29%  '

$ws.Range('E50').Value = '  +3.02%  '

$ws.Range('D51').Value = '52.24'
$ws.Range('E51').Value = '  +1.33%  '

# Restore the default "Normal" style so the cells keep their original
# (unstyled) formatting instead of the text-format style just applied.
$ws.Range("D2:E51").Style = "Normal"
